# New Submission Synced: 2026-02-07 19:35:18
#
# The "JSS 3B" responses sheet gained one new row from the form sync, and
# the previous row's "Admission No" (C2), which had been stored as text,
# got normalized to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3B")

# C2 was stored as text "19" - convert it to a real number (19).
$ws.Cells.Item(2, 3).Value = 19

# Append the new form-submission row (row 3).
$ws.Cells.Item(3, 1).Value = "2026-02-07 19:35:18"
$ws.Cells.Item(3, 2).Value = "Arhyel Jacob wakawa"

# Admission No keeps arriving as text from the form (unlike C2 above), e.g.
# "22". A plain Value assignment of a numeric-looking string would silently
# store it as a number, so build it as a TEXT() formula result in a scratch
# cell and paste-special just the value in - that keeps the destination
# cell's type as text without needing a "@" number-format override (which
# would permanently add an unused style record to the workbook).
$ws.Cells.Item(5, 5).Formula = "=TEXT(22,""0"")"
$ws.Cells.Item(5, 5).Copy()
$ws.Cells.Item(3, 3).PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Cells.Item(5, 5).Clear()

$ws.Cells.Item(3, 4).Value = 9
